# Applies the row-removal + re-randomization of the kitchens categorization
# stimulus order described in the commit "finished input list generation
# and sanity checks."
#
# One trial (stimuli/img_r77yy.png) is dropped from the list, the
# remaining 26 trials are reshuffled, and the trial_total (F) /
# trial_block (E) counters are renumbered to stay contiguous.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the very last data row so the sheet shrinks from 28 to 27 rows
# (A1:S28 -> A1:S27); the remaining rows 2-27 are fully overwritten below
# with the new trial order, so it does not matter which physical row is
# removed as long as the row count is reduced by one.
$ws.Rows(28).Delete()

# New contents for columns F,G,H,I,K,L,M,N,O,P,Q,R,S, in row order (row 2..27)
$data = @(
    @(217,"kitchens","kitchens","target","j","stimuli/img_aplao.png",64.0909090909091,40.75757575757576,52.42424242424242,33,3,3,3),
    @(218,"kitchens","kitchens","target","j","stimuli/img_yeh72.png",68.66666666666667,45.21212121212121,56.93939393939394,33,4,4,4),
    @(219,"kitchens","kitchens","target","j","stimuli/img_wyl6z.png",59.8235294117647,36.23529411764706,48.02941176470588,34,3,3,3),
    @(220,"kitchens","kitchens","target","j","stimuli/img_inqod.png",70.84848484848484,50.63636363636363,60.74242424242424,33,5,5,5),
    @(221,"kitchens","kitchens","target","j","stimuli/img_nyv2b.png",11.91176470588235,6.852941176470588,9.382352941176471,34,1,1,1),
    @(222,"kitchens","kitchens","target","j","stimuli/img_3gm8h.png",65.07894736842105,43.92105263157895,54.5,38,4,4,4),
    @(223,"kitchens","kitchens","target","j","stimuli/img_t90e2.png",83.0625,61.96875,72.515625,32,9,9,9),
    @(224,"kitchens","kitchens","target","j","stimuli/img_a8wvq.png",86.25925925925925,66.25925925925925,76.25925925925925,27,10,10,10),
    @(225,"kitchens","living_rooms","distractor","f","stimuli/img_jpjeg.png",90.90697674418605,74.3953488372093,82.65116279069767,43,10,10,10),
    @(226,"kitchens","kitchens","target","j","stimuli/img_iyxnj.png",75.30555555555556,54.33333333333334,64.81944444444444,36,6,6,6),
    @(227,"kitchens","living_rooms","distractor","f","stimuli/img_g13d5.png",73,51.51111111111111,62.25555555555556,45,6,6,6),
    @(228,"kitchens","living_rooms","distractor","f","stimuli/img_x9w7o.png",92.38888888888889,72.94444444444444,82.66666666666666,36,10,10,10),
    @(229,"kitchens","kitchens","target","j","stimuli/img_ce9vx.png",75.9090909090909,57.12121212121212,66.51515151515152,33,7,7,7),
    @(230,"kitchens","kitchens","target","j","stimuli/img_9mky8.png",84.32352941176471,65.17647058823529,74.75,34,9,9,9),
    @(231,"kitchens","kitchens","target","j","stimuli/img_uwv6y.png",78.88888888888889,59.30555555555556,69.09722222222223,36,8,8,8),
    @(232,"kitchens","kitchens","target","j","stimuli/img_60242.png",78.33333333333333,57.57575757575758,67.95454545454545,33,7,7,7),
    @(233,"kitchens","kitchens","target","j","stimuli/img_eatdk.png",81.40625,61.375,71.390625,32,8,8,8),
    @(234,"kitchens","kitchens","target","j","stimuli/img_7wul8.png",43.03030303030303,25.54545454545455,34.28787878787879,33,1,1,1),
    @(235,"kitchens","bedrooms","distractor","f","stimuli/img_scrdm.png",78.675,57.9,68.2875,40,7,7,7),
    @(236,"kitchens","kitchens","target","j","stimuli/img_ye5sl.png",53.2258064516129,34.45161290322581,43.83870967741936,31,2,2,2),
    @(237,"kitchens","kitchens","target","j","stimuli/img_d8xbu.png",91.36363636363636,73.18181818181819,82.27272727272728,33,10,10,10),
    @(238,"kitchens","kitchens","target","j","stimuli/img_463mq.png",51.35294117647059,30.20588235294118,40.77941176470588,34,2,2,2),
    @(239,"kitchens","living_rooms","distractor","f","stimuli/img_53nbn.png",73.28888888888889,51.15555555555556,62.22222222222223,45,6,6,6),
    @(240,"kitchens","kitchens","target","j","stimuli/img_p3hpc.png",72.83333333333333,52.22222222222222,62.52777777777777,36,6,6,6),
    @(241,"kitchens","bedrooms","distractor","f","stimuli/img_kn0we.png",80.1590909090909,56.68181818181818,68.42045454545455,44,7,7,7),
    @(242,"kitchens","kitchens","target","j","stimuli/img_cnyac.png",69.1470588235294,47.8235294117647,58.48529411764706,34,5,5,5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 6).Value  = $row[0]   # F trial_total
    $ws.Cells.Item($r, 7).Value  = $row[1]   # G target_cat
    $ws.Cells.Item($r, 8).Value  = $row[2]   # H category
    $ws.Cells.Item($r, 9).Value  = $row[3]   # I cond_cat
    $ws.Cells.Item($r, 11).Value = $row[4]   # K correct_answer
    $ws.Cells.Item($r, 12).Value = $row[5]   # L stimulus
    $ws.Cells.Item($r, 13).Value = $row[6]   # M conceptual
    $ws.Cells.Item($r, 14).Value = $row[7]   # N perceptual
    $ws.Cells.Item($r, 15).Value = $row[8]   # O typicality
    $ws.Cells.Item($r, 16).Value = $row[9]   # P n
    $ws.Cells.Item($r, 17).Value = $row[10]  # Q p_typicality
    $ws.Cells.Item($r, 18).Value = $row[11]  # R p_conceptual
    $ws.Cells.Item($r, 19).Value = $row[12]  # S p_perceptual
}
